$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 50002096
$ws.Range("J112").Value = 2551.1875
$ws.Range("L112").Value = 7653.5625
$ws.Range("N112").Value = -9869.5625
$ws.Range("H129").Value = 849.9143
$ws.Range("J129").Value = 1006.4706
$ws.Range("L129").Value = 3019.4118
$ws.Range("N129").Value = -13019.4118
$ws.Range("H137").Value = 2943803.8
$ws.Range("I137").Value = 5264839
$ws.Range("J137").Value = 3825.5334
$ws.Range("K137").Value = 15794517
$ws.Range("L137").Value = 11476.6002
$ws.Range("M137").Value = -15791967
$ws.Range("N137").Value = -16576.6002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13857.454
$ws.Range("I32").Value = 12772.718
$ws.Range("K32").Value = 12772.718
$ws.Range("M32").Value = -12485.718
$ws.Range("H61").Value = 45547136
$ws.Range("I61").Value = 66734804
$ws.Range("K61").Value = 66734804
$ws.Range("M61").Value = -66734592
$ws.Range("H74").Value = 7638838.5
$ws.Range("I74").Value = 11953899
$ws.Range("J74").Value = 87483.336
$ws.Range("K74").Value = 11953899
$ws.Range("L74").Value = 87483.336
$ws.Range("M74").Value = -11953025
$ws.Range("N74").Value = -89231.336
$ws.Range("H77").Value = 7638838.5
$ws.Range("I77").Value = 11953899
$ws.Range("J77").Value = 87483.336
$ws.Range("K77").Value = 59769495
$ws.Range("L77").Value = 437416.68
$ws.Range("M77").Value = -59765127
$ws.Range("N77").Value = -446152.68
$ws.Range("H132").Value = 41162.727
$ws.Range("I132").Value = 32806.72
$ws.Range("J132").Value = 55236
$ws.Range("K132").Value = 98420.16
$ws.Range("L132").Value = 165708
$ws.Range("M132").Value = -95890.16
$ws.Range("N132").Value = -170768
$ws.Range("H136").Value = 45547136
$ws.Range("I136").Value = 66734804
$ws.Range("K136").Value = 200204412
$ws.Range("M136").Value = -200201862

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 808.61536
$ws.Range("I64").Value = 500.66666
$ws.Range("J64").Value = 901
$ws.Range("K64").Value = 500.66666
$ws.Range("L64").Value = 901
$ws.Range("M64").Value = -275.66666
$ws.Range("N64").Value = -1351
$ws.Range("H67").Value = 808.61536
$ws.Range("I67").Value = 500.66666
$ws.Range("J67").Value = 901
$ws.Range("K67").Value = 500.66666
$ws.Range("L67").Value = 901
$ws.Range("M67").Value = 279.33334
$ws.Range("N67").Value = -2461
$ws.Range("H134").Value = 2366.4167
$ws.Range("I134").Value = 1897.1923
$ws.Range("K134").Value = 5691.5769
$ws.Range("M134").Value = -3156.5769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1763.6945
$ws.Range("I31").Value = 1232.6333
$ws.Range("J31").Value = 4419
$ws.Range("K31").Value = 1232.6333
$ws.Range("L31").Value = 4419
$ws.Range("M31").Value = -937.6333
$ws.Range("N31").Value = -5009
$ws.Range("H34").Value = 1763.6945
$ws.Range("I34").Value = 1232.6333
$ws.Range("J34").Value = 4419
$ws.Range("K34").Value = 1232.6333
$ws.Range("L34").Value = 4419
$ws.Range("M34").Value = -1030.6333
$ws.Range("N34").Value = -4823
$ws.Range("H58").Value = 23257574
$ws.Range("I58").Value = 35715480
$ws.Range("K58").Value = 35715480
$ws.Range("M58").Value = -35715277
$ws.Range("H68").Value = 22314.166
$ws.Range("J68").Value = 22314.166
$ws.Range("L68").Value = 22314.166
$ws.Range("N68").Value = -23812.166
$ws.Range("H71").Value = 22314.166
$ws.Range("J71").Value = 22314.166
$ws.Range("L71").Value = 66942.49800000001
$ws.Range("N71").Value = -74430.49800000001
$ws.Range("H132").Value = 52881.6
$ws.Range("I132").Value = 2189
$ws.Range("J132").Value = 86676.664
$ws.Range("K132").Value = 6567
$ws.Range("L132").Value = 260029.992
$ws.Range("M132").Value = -4037
$ws.Range("N132").Value = -265089.992
$ws.Range("H134").Value = 29284.951
$ws.Range("I134").Value = 1846.6364
$ws.Range("K134").Value = 5539.9092
$ws.Range("M134").Value = -3004.9092
$ws.Range("H136").Value = 23257574
$ws.Range("I136").Value = 35715480
$ws.Range("K136").Value = 107146440
$ws.Range("M136").Value = -107143890

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 738.36365
$ws.Range("I5").Value = 564.3333
$ws.Range("J5").Value = 858.8461
$ws.Range("K5").Value = 1692.9999
$ws.Range("L5").Value = 2576.5383
$ws.Range("M5").Value = -1580.9999
$ws.Range("N5").Value = -2800.5383
$ws.Range("H122").Value = 884.5
$ws.Range("I122").Value = 386.8889
$ws.Range("K122").Value = 3482.0001
$ws.Range("M122").Value = -1032.0001
$ws.Range("H131").Value = 20899.021
$ws.Range("I131").Value = 580.8333
$ws.Range("J131").Value = 28287.455
$ws.Range("K131").Value = 1742.4999
$ws.Range("L131").Value = 84862.36500000001
$ws.Range("M131").Value = 3297.5001
$ws.Range("N131").Value = -94942.36500000001
$ws.Range("H132").Value = 1329.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 1329.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 11965.5
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -17025.5
$ws.Range("H135").Value = 738.36365
$ws.Range("I135").Value = 564.3333
$ws.Range("J135").Value = 858.8461
$ws.Range("K135").Value = 5078.9997
$ws.Range("L135").Value = 7729.6149
$ws.Range("M135").Value = -2543.9997
$ws.Range("N135").Value = -12799.6149
$ws.Range("H136").Value = 1400
$ws.Range("I136").Value = 1400
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4200
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 900
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 30301.205
$ws.Range("I70").Value = 41387
$ws.Range("J70").Value = 5358.1665
$ws.Range("K70").Value = 41387
$ws.Range("L70").Value = 5358.1665
$ws.Range("M70").Value = -41117
$ws.Range("N70").Value = -5898.1665
$ws.Range("H73").Value = 30301.205
$ws.Range("I73").Value = 41387
$ws.Range("J73").Value = 5358.1665
$ws.Range("K73").Value = 41387
$ws.Range("L73").Value = 5358.1665
$ws.Range("M73").Value = -40451
$ws.Range("N73").Value = -7230.1665
$ws.Range("H132").Value = 68401.7
$ws.Range("I132").Value = 63534.375
$ws.Range("J132").Value = 73964.36
$ws.Range("K132").Value = 190603.125
$ws.Range("L132").Value = 221893.08
$ws.Range("M132").Value = -188073.125
$ws.Range("N132").Value = -226953.08

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 20934.096
$ws.Range("I132").Value = 1404.2683
$ws.Range("K132").Value = 4212.8049
$ws.Range("M132").Value = -1682.8049
$ws.Range("H136").Value = 96604.38
$ws.Range("I136").Value = 63661.75
$ws.Range("J136").Value = 202020.8
$ws.Range("K136").Value = 190985.25
$ws.Range("L136").Value = 606062.3999999999
$ws.Range("M136").Value = -188435.25
$ws.Range("N136").Value = -611162.3999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 46608.66
$ws.Range("I132").Value = 30085.97
$ws.Range("J132").Value = 102785.8
$ws.Range("K132").Value = 90257.91
$ws.Range("L132").Value = 308357.4
$ws.Range("M132").Value = -87727.91
$ws.Range("N132").Value = -313417.4
$ws.Range("H136").Value = 49044.594
$ws.Range("I136").Value = 28863
$ws.Range("J136").Value = 170134.17
$ws.Range("K136").Value = 86589
$ws.Range("L136").Value = 510402.51
$ws.Range("M136").Value = -84039
$ws.Range("N136").Value = -515502.51
